$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Regenerated s_val data (filtering save games) for rows 2-5, columns B:E and G.
# Column A (date) and column F (Win) are unchanged.

$ws.Range("B2").Value = 0.7287194209349384
$ws.Range("C2").Value = 1.65323645889881
$ws.Range("D2").Value = 0.1529057820181812
$ws.Range("E2").Value = 0.4998867070740569
$ws.Range("G2").Value = 3.034748368925986

$ws.Range("B3").Value = 3.182878228561681
$ws.Range("C3").Value = 1.65323645889881
$ws.Range("D3").Value = 3.082599426703578
$ws.Range("E3").Value = 0.4998867070740569
$ws.Range("G3").Value = 8.418600821238126

$ws.Range("B4").Value = 3.182878228561681
$ws.Range("C4").Value = 0.3375848360084654
$ws.Range("D4").Value = 0.7127328510149897
$ws.Range("E4").Value = 0.4998867070740569
$ws.Range("G4").Value = 4.733082622659194

$ws.Range("B5").Value = 0.02258322285507441
$ws.Range("C5").Value = 0.004309184025731883
$ws.Range("D5").Value = 16.98373111632243
$ws.Range("E5").Value = 246.9852506941017
$ws.Range("G5").Value = 263.9958742173049
